{"js": "// Update the worksheet date line and the 25 division problems in the\n// table to the new values. Every populated paragraph in the document\n// changes (the title date plus each filled-in table cell). We address\n// each table cell positionally by (row, column) rather than searching\n// for the old text, because some of the old division expressions (e.g.\n// \"57\u00f79=\") repeat verbatim elsewhere in the document with a different\n// replacement, so a text-based search/replace would be ambiguous.\n//\n// Each cell/paragraph already contains a single run with the\n// formatting (font, size, alignment) that must be kept, so we replace\n// the paragraph's text in place (InsertLocation.replace) instead of\n// rewriting the whole cell body, which would drop that formatting.\n\nconst body = context.document.body;\n\n// Title / date line (first paragraph in the document body).\nconst titleParagraphs = body.paragraphs;\ntitleParagraphs.load(\"items\");\n\n// The worksheet table: only rows 0, 4, 8, 12, 16 (of 20, 0-indexed)\n// hold the division problems; the rows between them are blank answer\n// rows.\nconst tables = body.tables;\ntables.load(\"items\");\n\nawait context.sync();\n\ntitleParagraphs.items[0].insertText(\"2023-12-11 Monday\", Word.InsertLocation.replace);\n\nconst table = tables.items[0];\nconst newValues = {\n  0: [\"77\u00f78=\", \"38\u00f76=\", \"42\u00f73=\", \"33\u00f73=\", \"62\u00f75=\"],\n  4: [\"15\u00f75=\", \"30\u00f72=\", \"62\u00f74=\", \"68\u00f77=\", \"85\u00f72=\"],\n  8: [\"14\u00f74=\", \"57\u00f73=\", \"15\u00f77=\", \"71\u00f77=\", \"85\u00f79=\"],\n  12: [\"93\u00f78=\", \"52\u00f73=\", \"88\u00f74=\", \"96\u00f76=\", \"67\u00f73=\"],\n  16: [\"56\u00f73=\", \"19\u00f79=\", \"69\u00f78=\", \"98\u00f72=\", \"37\u00f76=\"],\n};\n\nconst cellParagraphs = [];\nfor (const [rowStr, rowValues] of Object.entries(newValues)) {\n  const row = parseInt(rowStr, 10);\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.body.paragraphs.load(\"items\");\n    cellParagraphs.push({ paragraphs: cell.body.paragraphs, text: rowValues[col] });\n  }\n}\n\nawait context.sync();\n\nfor (const { paragraphs, text } of cellParagraphs) {\n  paragraphs.items[0].insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date line and the 25 division problems in the\n# table to the new values. Every populated paragraph in the document\n# changes (the title date plus each filled-in table cell). We address\n# each table cell positionally by (row, column) rather than searching\n# for the old text, because some of the old division expressions (e.g.\n# \"57\u00f79=\") repeat verbatim elsewhere in the document with a different\n# replacement, so a text-based Find/Replace would be ambiguous.\n\n$d = $word.ActiveDocument\n\n# Title / date line (first paragraph in the document body).\n$d.Paragraphs(1).Range.Text = \"2023-12-11 Monday\"\n\n# The worksheet table: only rows 1, 5, 9, 13, 17 (of 20) hold the\n# division problems; the rows between them are blank answer rows.\n$t = $d.Tables(1)\n\n$newValues = @{\n    1  = @(\"77\u00f78=\", \"38\u00f76=\", \"42\u00f73=\", \"33\u00f73=\", \"62\u00f75=\")\n    5  = @(\"15\u00f75=\", \"30\u00f72=\", \"62\u00f74=\", \"68\u00f77=\", \"85\u00f72=\")\n    9  = @(\"14\u00f74=\", \"57\u00f73=\", \"15\u00f77=\", \"71\u00f77=\", \"85\u00f79=\")\n    13 = @(\"93\u00f78=\", \"52\u00f73=\", \"88\u00f74=\", \"96\u00f76=\", \"67\u00f73=\")\n    17 = @(\"56\u00f73=\", \"19\u00f79=\", \"69\u00f78=\", \"98\u00f72=\", \"37\u00f76=\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($col = 1; $col -le $rowValues.Count; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
